# [Fonds de solidarite] Add 2020-08-18 data
#
# Refreshes "nombre_aides" (col C) and "montant_total" (col D) for the
# (region x classe_effectif) rows impacted by the 2020-08-18 data pull.
# Every value in this sheet is stored as text (the source export writes
# inline strings, not numbers), so each target cell is switched to the
# Text number format before the write — this keeps the literal string
# representation (e.g. the trailing ".00") instead of Excel coercing the
# input into a numeric value and dropping formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2;   C = "182";  D = "422816.00" },
    @{ Row = 3;   C = "979";  D = "3084177.33" },
    @{ Row = 4;   C = "407";  D = "1639698.25" },
    @{ Row = 5;   C = "112";  D = "527128.09" },
    @{ Row = 6;   C = "28";   D = "183143.82" },
    @{ Row = 7;   C = "7";    D = "40500.00" },
    @{ Row = 8;   C = "35";   D = "70000.00" },
    @{ Row = 17;  C = "151";  D = "642202.10" },
    @{ Row = 33;  C = "105";  D = "285673.00" },
    @{ Row = 34;  C = "556";  D = "1791904.47" },
    @{ Row = 35;  C = "222";  D = "1108288.11" },
    @{ Row = 47;  C = "44";   D = "256060.00" },
    @{ Row = 50;  C = "98";   D = "275768.17" },
    @{ Row = 51;  C = "564";  D = "1929686.52" },
    @{ Row = 52;  C = "258";  D = "1116446.76" },
    @{ Row = 53;  C = "87";   D = "509378.23" },
    @{ Row = 80;  C = "877";  D = "2787444.67" },
    @{ Row = 81;  C = "334";  D = "1341403.79" },
    @{ Row = 83;  C = "28";   D = "164080.04" },
    @{ Row = 97;  C = "293";  D = "769279.43" },
    @{ Row = 98;  C = "1201"; D = "3655891.37" },
    @{ Row = 99;  C = "448";  D = "1816594.02" },
    @{ Row = 100; C = "119";  D = "544000.00" },
    @{ Row = 101; C = "32";   D = "199157.00" }
)

foreach ($chg in $changes) {
    $r = $chg.Row
    $rowRange = $ws.Range("C" + $r + ":D" + $r)
    $rowRange.NumberFormat = "@"
    $ws.Range("C" + $r).Value = $chg.C
    $ws.Range("D" + $r).Value = $chg.D
}
